$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2) currently holds a mailto hyperlink styled with the built-in
# "Hyperlink" cell style; the target just has plain text there.
$ws.Range("A2").ClearFormats()
$ws.Hyperlinks.Delete()

# Add the two new rows of credentials.
$ws.Range("A3").Value = "tzw@gmail.com"
$ws.Range("B3").Value = "tzw123"

$ws.Range("A4").Value = "example@gmail.com"

# Force B4 to be stored as text (not a number) by setting the number
# format to Text before assigning the value, then clear the leftover
# number-format so the cell keeps the default (unstyled) appearance.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "123"
$ws.Range("B4").ClearFormats()
